# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest scrape (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (row -> new value for column F)
$zhanlanUpdates = @{
    2  = 2245
    9  = 2641
    11 = 1656
    13 = 273
    14 = 678
    16 = 123
    17 = 346
    20 = 43
    22 = 5836
    23 = 238
    24 = 1086
    25 = 124
    26 = 168
    27 = 148
    28 = 270
    29 = 237
    31 = 1068
    36 = 431
    37 = 1218
    42 = 137
}

$wsZhanlan = $wb.Worksheets.Item("展览")
foreach ($row in $zhanlanUpdates.Keys) {
    $wsZhanlan.Cells.Item($row, 6).Value = $zhanlanUpdates[$row]
}

# Sheet: 全部类型 (row -> new value for column F)
$allTypesUpdates = @{
    2  = 2245
    14 = 2641
    16 = 1656
    18 = 273
    19 = 678
    22 = 123
    23 = 346
    25 = 43
    27 = 5836
    28 = 238
    29 = 1086
    30 = 124
    31 = 168
    32 = 148
    33 = 270
    34 = 237
    36 = 1068
    40 = 431
    41 = 1218
    46 = 137
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
